$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends at column R (2021). Add a new column S for 2022,
# mirroring the formatting of column R (the previous "year" column).

# Copy the formatting of R2:R5 onto the new S2:S5 cells first, so the new
# column inherits borders/fonts/number formats exactly like the rest of the
# table (this is what a human would do in Excel: copy the last column and
# paste-special its formats into the new one).
$ws.Range("R2:R5").Copy()
$ws.Range("S2:S5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 is an empty bottom-border spacer row - no value needed for S2.

# Row 3: year header.
$ws.Range("S3").Value = 2022

# Row 4: population count.
$ws.Range("S4").Value = 211650

# Row 5: percentage of total population (same published value as 2021).
$ws.Range("S5").Value = 2.9794303052841493

# Match the author's selection on the newly added column.
$ws.Range("S2").Select()
